$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 380, shifting the existing rows 380-406 down to 381-407.
$ws.Rows.Item(380).EntireRow.Insert()

# Populate the newly inserted row 380 with the new weekly price record.
$ws.Cells.Item(380, 1).Value  = 4
$ws.Cells.Item(380, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(380, 3).Value  = "Los Lagos"
$ws.Cells.Item(380, 4).Value  = 44826
$ws.Cells.Item(380, 5).Value  = 10
$ws.Cells.Item(380, 6).Value  = 100112023
$ws.Cells.Item(380, 7).Value  = "Brócoli"
$ws.Cells.Item(380, 8).Value  = "Sin especificar"
$ws.Cells.Item(380, 9).Value  = "Primera"
$ws.Cells.Item(380, 10).Value = 500
$ws.Cells.Item(380, 11).Value = 1500
$ws.Cells.Item(380, 12).Value = 1500
$ws.Cells.Item(380, 13).Value = 1500
$ws.Cells.Item(380, 14).Value = "$/unidad"
$ws.Cells.Item(380, 15).Value = "Región Metropolitana"
$ws.Cells.Item(380, 16).Value = 1500
$ws.Cells.Item(380, 17).Value = 1
$ws.Cells.Item(380, 18).Value = "Hortaliza"
